$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E22").Value = "Delincuencia"
$ws.Range("E23").Value = "Educación"
$ws.Range("E24").Value = "Salud"
$ws.Range("E25").Value = "Empleo"
$ws.Range("E26").Value = "Pensiones"
$ws.Range("E27").Value = "Medio ambiente"
$ws.Range("E28").Value = "Inmigración"
$ws.Range("E29").Value = "Derechos para la mujer"
$ws.Range("E30").Value = "Economía"
